# Remove the post entry for "「私できるよ」" which occupied row 820.
# Deleting the entire row shifts all subsequent rows (821-836) up by one,
# so they become rows 820-835, and the used range shrinks from
# A1:C836 to A1:C835.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(820).Delete()
